# KIM_YR_FIN.xlsx - add newest fiscal-year column (FY ending 2018-12-31)
# This inserts a new column D (shifting old D:K -> E:L) on the single
# worksheet, carries the existing column formatting into the new column,
# and fills in the new year's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before column D; old D:K becomes E:L.
$ws.Columns("D").Insert()

# 2) Copy number formats/styles from the (new) column E into the
#    newly inserted column D, restricted to the three data blocks so we
#    do not materialise cells in rows that have none (16/36/78 gaps etc.)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# 3) Match the column width of the newly inserted column to its neighbors.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# 4) Populate column D with the new fiscal year's values.
$values = @{
    "7" = 43465
    "8" = 1164800
    "9" = 328600
    "10" = 836200
    "12" = "NA"
    "13" = 0
    "14" = 92000
    "15" = 310400
    "17" = 825000
    "18" = 339800
    "20" = 13000
    "21" = 663200
    "22" = 183300
    "23" = 169500
    "24" = 1600
    "25" = 0
    "26" = 167900
    "27" = 199900
    "28" = 0
    "29" = 229800
    "30" = 0
    "31" = 0
    "32" = -13000
    "33" = 429700
    "34" = 0
    "35" = 429700
    "38" = 43465
    "41" = 143600
    "42" = 0
    "43" = 184500
    "44" = 0
    "45" = 0
    "46" = 0
    "47" = 787800
    "48" = 9491900
    "49" = 0
    "50" = 0
    "51" = 0
    "52" = 156200
    "53" = 0
    "54" = 10999100
    "57" = 174900
    "58" = 0
    "59" = 130300
    "60" = 0
    "61" = 4873900
    "62" = "NA"
    "63" = 0
    "64" = 0
    "65" = 0
    "66" = 5665300
    "68" = 0
    "69" = 0
    "70" = 0
    "71" = 0
    "72" = -787700
    "73" = 0
    "74" = 0
    "75" = 0
    "76" = 5333800
    "77" = 0
    "80" = 43465
    "81" = 429700
    "83" = 310400
    "84" = 0
    "85" = 0
    "86" = 0
    "87" = 0
    "88" = 0
    "89" = 637900
    "91" = -536900
    "92" = 0
    "93" = 0
    "94" = 253600
    "96" = -529800
    "97" = 0
    "98" = 0
    "99" = 0
    "100" = -986500
    "101" = 0
    "102" = -94900
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
